$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.436.49"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.654.92"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.44"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.13"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.541"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.653.55"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.06"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.138.24"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.366.86"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.622.65"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  +3.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "364.61"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.49"
$ws.Range("E22").Value = "  +2.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.84"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "75.37"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  +4.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.783.12"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000103"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "576.64"
$ws.Range("E31").Value = "  +3.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("E32").Value = "  +3.17%  "
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.65"
$ws.Range("E35").Value = "  +5.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.130"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.88"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.71"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.374"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0317"
$ws.Range("E45").Value = "  -7.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.79"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.94"
$ws.Range("E47").Value = "  +5.47%  "
$ws.Range("E48").Value = "  +3.01%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.595"
$ws.Range("E49").Value = "  +6.75%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.97"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("E51").Value = "  -0.16%  "
